$wb = $excel.ActiveWorkbook

# Helper that writes a block of rows (each an array of 6 strings: Date,
# Timestamp, Hour, Location, Value, Status) into a worksheet starting at
# $startRow. Column A always holds a literal "YYYY-MM-DD" text value, so it
# is apostrophe-prefixed to stop Excel from converting it into a date serial
# number. Column E is only apostrophe-prefixed when it looks like a bare
# percentage (e.g. "88.5%"), which Excel would otherwise reinterpret as a
# numeric percentage.
function Fill-Rows {
    param(
        $ws,
        [int]$startRow,
        $rowsData
    )
    for ($i = 0; $i -lt $rowsData.Count; $i++) {
        $r = $startRow + $i
        $rd = $rowsData[$i]
        $ws.Cells.Item($r, 1).Value = "'" + $rd[0]
        $ws.Cells.Item($r, 2).Value = $rd[1]
        $ws.Cells.Item($r, 3).Value = $rd[2]
        $ws.Cells.Item($r, 4).Value = $rd[3]
        if ($rd[4] -match "^\d+(\.\d+)?%$") {
            $ws.Cells.Item($r, 5).Value = "'" + $rd[4]
        } else {
            $ws.Cells.Item($r, 5).Value = $rd[4]
        }
        $ws.Cells.Item($r, 6).Value = $rd[5]
    }
}

# --- ALERTS sheet: add rows 2-7 ---
$wsAlerts = $wb.Worksheets.Item("ALERTS")
$alertsData = @(
    @("2026-01-30", "14:57:17", "14:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "14:57:20", "14:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "14:59:04", "14:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:00:46", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:00:49", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:03:04", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED")
)
Fill-Rows $wsAlerts 2 $alertsData

# --- mmWave sheet: add rows 27-53 ---
$wsMmWave = $wb.Worksheets.Item("mmWave")
$mmWaveData = @(
    @("2026-01-30", "14:55:23", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:55:33", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:55:44", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:56:10", "14:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "14:57:17", "14:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "14:57:20", "14:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "14:57:24", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:57:34", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:57:44", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:57:55", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:58:06", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:58:16", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:58:26", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:58:37", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:59:04", "14:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "14:59:50", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:00:00", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:00:10", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:00:46", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:00:49", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:02:23", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:02:34", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:02:44", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:03:04", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:03:51", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:03:54", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:03:58", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED")
)
Fill-Rows $wsMmWave 27 $mmWaveData

# --- PIR sheet: add rows 14-21 ---
$wsPir = $wb.Worksheets.Item("PIR")
$pirData = @(
    @("2026-01-30", "14:57:24", "14:00", "Living Room", "RECOVERY_DETECTION", "Inactive"),
    @("2026-01-30", "14:59:50", "14:00", "Living Room", "RECOVERY_DETECTION", "Inactive"),
    @("2026-01-30", "15:02:12", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:02:17", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:02:22", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:02:23", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive"),
    @("2026-01-30", "15:02:27", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:02:32", "15:00", "Bathroom", "No Motion", "Inactive")
)
Fill-Rows $wsPir 14 $pirData

# --- Humidity sheet: add rows 13-17 ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityData = @(
    @("2026-01-30", "15:02:12", "15:00", "Bathroom", "88.5%", "Active"),
    @("2026-01-30", "15:02:17", "15:00", "Bathroom", "88.5%", "Active"),
    @("2026-01-30", "15:02:22", "15:00", "Bathroom", "88.4%", "Active"),
    @("2026-01-30", "15:02:27", "15:00", "Bathroom", "87.5%", "Active"),
    @("2026-01-30", "15:02:32", "15:00", "Bathroom", "88.5%", "Active")
)
Fill-Rows $wsHumidity 13 $humidityData
